$wb = $excel.ActiveWorkbook

# Rename sheet tabs (new timestamp-based names)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477869264699"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778711436625"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477871149665"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477871208665"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650477871271667"

# Sheet 1 (GNG_TO) - column B values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778692226634.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778692477016.csv"
$ws1.Range("B4").Value = "go_stims-16504778692486663.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477869263701.csv"

# Sheet 2 (NB_TO) - column B values
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504778703346655.csv"
$ws2.Range("B3").Value = "OB-16504778702466648.csv"
$ws2.Range("B4").Value = "ZB-match_9-16504778696316636.csv"
$ws2.Range("B5").Value = "TB-16504778709887002.csv"
$ws2.Range("B6").Value = "OB-16504778703186638.csv"
$ws2.Range("B7").Value = "TB-16504778704427018.csv"
$ws2.Range("B8").Value = "TB-16504778711206942.csv"
$ws2.Range("B9").Value = "ZB-match_8-16504778698396957.csv"
$ws2.Range("B10").Value = "ZB-match_5-16504778699026983.csv"

# Sheet 3 (RS_TO) - column B values (swap eyes closed/open)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL_TO) - column B values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778711757002.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778711506636.csv"
$ws4.Range("B4").Value = "MM_stims-16504778711916652.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778711767004.csv"
$ws4.Range("B6").Value = "MM_stims-1650477871207698.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778711926649.csv"

# Sheet 5 (vSAT_TO) - column B values
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778712556987.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778712116694.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778712236633.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778712397.csv"
